$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.362.97"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "1.667.61"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.78"
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5168"
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06454"
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2566"
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.99"
$ws.Range("E10").Value = "  -4.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07661"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.671.76"
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.323"
$ws.Range("E13").Value = "  -5.45%  "
$ws.Range("D14").Value = "1.895.85"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5538"
$ws.Range("E15").Value = "  -3.79%  "
$ws.Range("D16").Value = "0.0₅8049"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.55"
$ws.Range("E17").Value = "  -4.86%  "
$ws.Range("D18").Value = "26.389.48"
$ws.Range("E18").Value = "  -3.53%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.82"
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.405"
$ws.Range("E21").Value = "  -5.74%  "
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.900"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.76"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.755"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1164"
$ws.Range("E27").Value = "  -4.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.003"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("E29").Value = "  -3.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05268"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.373"
$ws.Range("E32").Value = "  -3.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.222"
$ws.Range("E33").Value = "  -6.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.569"
$ws.Range("E34").Value = "  -4.71%  "
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.378"
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9287"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5733"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("D39").Value = "1.154.88"
$ws.Range("E39").Value = "  +10.64%  "
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8474"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.008"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.661"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.01"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "1.805.50"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -5.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4501"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.09"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.909"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("E51").Value = "  -2.57%  "
